$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 08:22"

# --- Row 4: Estados Unidos (totals update, name unchanged) ---
$ws.Range("B4").Value = 925758
$ws.Range("C4").Value = 526
$ws.Range("D4").Value = 110432
$ws.Range("E4").Value = 763109
$ws.Range("F4").Value = 15097
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 52217

# --- Rows 64/65: Hungria/Kazajistan swap order & update data ---
# Row 64 becomes Kazajistan with fresh numbers; row 65 becomes Hungria
# carrying what used to be row 64's (Hungria) figures.
$ws.Range("A64").Value = "Kazajistan"
$ws.Range("B64").Value = 2482
$ws.Range("C64").Value = 66
$ws.Range("D64").Value = 604
$ws.Range("E64").Value = 1853
$ws.Range("F64").Value = 29
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 25

$ws.Range("A65").Value = "Hungria"
$ws.Range("B65").Value = 2443
$ws.Range("C65").Value = 60
$ws.Range("D65").Value = 458
$ws.Range("E65").Value = 1723
$ws.Range("F65").Value = 61
$ws.Range("G65").Value = 12
$ws.Range("H65").Value = 262

# --- Rows 125/126: Vietnam/El Salvador swap order & update data ---
# Row 125 becomes El Salvador with fresh numbers; row 126 becomes Vietnam
# carrying what used to be row 125's (Vietnam) figures.
$ws.Range("A125").Value = "El Salvador"
$ws.Range("B125").Value = 274
$ws.Range("C125").Value = 13
$ws.Range("D125").Value = 75
$ws.Range("E125").Value = 191
$ws.Range("F125").Value = 2
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 8

$ws.Range("A126").Value = "Vietnam"
$ws.Range("B126").Value = 270
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 220
$ws.Range("E126").Value = 50
$ws.Range("F126").Value = 8
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

# --- Row 196: Montserrat (totals update, name unchanged) ---
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = 1
$ws.Range("G196").Value = 1
$ws.Range("H196").Value = 1
